$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cellAddr, $value)
    $rng = $ws.Range($cellAddr)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = $origStyle
}

Set-TextValue "D2" "27.170.54"
Set-TextValue "E2" "  +0.04%  "
Set-TextValue "D3" "1.827.83"
Set-TextValue "E3" "  -0.23%  "
Set-TextValue "D4" "1.011"
Set-TextValue "E4" "  +0.20%  "
Set-TextValue "D5" "313.47"
Set-TextValue "E5" "  +0.16%  "
Set-TextValue "D6" "1.012"
Set-TextValue "E6" "  +0.37%  "
Set-TextValue "D7" "0.4645"
Set-TextValue "E7" "  -1.25%  "
Set-TextValue "D8" "0.3645"
Set-TextValue "E8" "  -1.23%  "
Set-TextValue "D9" "0.07328"
Set-TextValue "E9" "  -1.07%  "
Set-TextValue "D10" "0.8764"
Set-TextValue "E10" "  -0.56%  "
Set-TextValue "D11" "20.23"
Set-TextValue "E11" "  -1.00%  "
Set-TextValue "D12" "1.884.47"
Set-TextValue "E12" "  +2.77%  "
Set-TextValue "D13" "0.07652"
Set-TextValue "E13" "  +4.24%  "
Set-TextValue "D14" "5.358"
Set-TextValue "E14" "  -2.10%  "
Set-TextValue "D15" "92.82"
Set-TextValue "E15" "  +0.01%  "
Set-TextValue "D16" "6.489"
Set-TextValue "E16" "  -0.97%  "
Set-TextValue "D17" "1.009"
Set-TextValue "E17" "  +0.01%  "
Set-TextValue "D18" "0.000008671"
Set-TextValue "E18" "  -1.28%  "
Set-TextValue "D19" "1.011"
Set-TextValue "E19" "  +0.36%  "
Set-TextValue "D20" "27.562.22"
Set-TextValue "E20" "  +1.39%  "
Set-TextValue "D21" "14.55"
Set-TextValue "E21" "  -1.62%  "
Set-TextValue "D22" "5.230"
Set-TextValue "E22" "  -1.43%  "
Set-TextValue "D23" "10.59"
Set-TextValue "E23" "  -0.86%  "
Set-TextValue "D24" "2.098.01"
Set-TextValue "E24" "  +1.82%  "
Set-TextValue "D25" "1.884"
Set-TextValue "E25" "  -0.77%  "
Set-TextValue "E26" "  -0.61%  "
Set-TextValue "E27" "  -0.65%  "
Set-TextValue "D28" "2.094"
Set-TextValue "E28" "  -2.90%  "
Set-TextValue "D29" "5.130"
Set-TextValue "E29" "  -2.66%  "
Set-TextValue "D30" "116.46"
Set-TextValue "E30" "  -0.84%  "
Set-TextValue "D31" "0.08919"
Set-TextValue "E31" "  -0.09%  "
Set-TextValue "D32" "2.962"
Set-TextValue "E32" "  +0.94%  "
Set-TextValue "B33" "ImmutableX"
Set-TextValue "C33" "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue "D33" "0.7382"
Set-TextValue "E33" "  -2.93%  "
Set-TextValue "B34" "ARBITRUM"
Set-TextValue "C34" "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue "D34" "1.155"
Set-TextValue "E34" "  -1.28%  "
Set-TextValue "B35" "Filecoin"
Set-TextValue "C35" "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue "D35" "4.477"
Set-TextValue "E35" "  -1.46%  "
Set-TextValue "E36" "  +0.22%  "
Set-TextValue "D37" "2.547"
Set-TextValue "E37" "  +5.49%  "
Set-TextValue "D38" "1.088"
Set-TextValue "E38" "  -1.30%  "
Set-TextValue "D39" "0.05260"
Set-TextValue "E39" "  -1.42%  "
Set-TextValue "D40" "0.01924"
Set-TextValue "E40" "  -1.86%  "
Set-TextValue "B41" "FraxShare"
Set-TextValue "C41" "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue "D41" "7.276"
Set-TextValue "E41" "  -0.61%  "
Set-TextValue "B42" "MXToken"
Set-TextValue "C42" "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue "D42" "2.933"
Set-TextValue "E42" "  -2.16%  "
Set-TextValue "D43" "0.5246"
Set-TextValue "E43" "  -1.86%  "
Set-TextValue "D44" "0.1634"
Set-TextValue "E44" "  -1.80%  "
Set-TextValue "D45" "8.322"
Set-TextValue "E45" "  -2.74%  "
Set-TextValue "D46" "0.4866"
Set-TextValue "E46" "  -1.58%  "
Set-TextValue "B47" "PaxDollar"
Set-TextValue "C47" "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-TextValue "D47" "1.013"
Set-TextValue "E47" "  +0.44%  "
Set-TextValue "B48" "EnergySwap"
Set-TextValue "C48" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D48" "10.27"
Set-TextValue "E48" "  -2.64%  "
Set-TextValue "D49" "103.64"
Set-TextValue "E49" "  -0.19%  "
Set-TextValue "D50" "1.640"
Set-TextValue "E50" "  -1.81%  "
Set-TextValue "D51" "0.06278"
Set-TextValue "E51" "  -0.67%  "
